$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Update the four changed values
$ws.Range("B7").Value = "aaadminunverified"
$ws.Range("B24").Value = "testpass3"
$ws.Range("B30").Value = "userone1"
$ws.Range("B34").Value = 1234567890

# Update the selection shown in the saved view
$ws.Activate()
$ws.Range("D6").Select()
